$wb = $excel.ActiveWorkbook

$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsTasas = $wb.Worksheets.Item("tasas")

# Update the "Conversión del día" message on Hoja1!A1 with the new rates
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.25 = 29348.04 pesos`n✅ 29348.04 pesos = 7.22 = 924.58 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Update the numeric rate values on the tasas sheet
$wsTasas.Range("N10").Value = 137.999
$wsTasas.Range("O10").Value = 4050
$wsTasas.Range("N12").Value = 4063
$wsTasas.Range("O12").Value = 128
